$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H2").Value = 1043.9286
$ws.Range("I2").Value = 1198.8889
$ws.Range("J2").Value = 765
$ws.Range("K2").Value = 1198.8889
$ws.Range("L2").Value = 765
$ws.Range("M2").Value = -1085.8889
$ws.Range("N2").Value = -991

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H9").Value = 141
$ws.Range("I9").Value = 143.33333
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 143.33333
$ws.Range("L9").Value = 120
$ws.Range("M9").Value = 25.66667000000001
$ws.Range("N9").Value = -458

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 395.54544
$ws.Range("I2").Value = 381.93332
$ws.Range("J2").Value = 424.7143
$ws.Range("K2").Value = 381.93332
$ws.Range("L2").Value = 424.7143
$ws.Range("M2").Value = -268.93332
$ws.Range("N2").Value = -650.7143

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H37").Value = 22768.625
$ws.Range("I37").Value = 2000
$ws.Range("K37").Value = 2000
$ws.Range("M37").Value = -1727

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H44").Value = 27840.834
$ws.Range("J44").Value = 30009
$ws.Range("L44").Value = 30009
$ws.Range("N44").Value = -30985

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H55").Value = 26999.2

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H63").Value = 31677.594
$ws.Range("I63").Value = 101995.336
$ws.Range("K63").Value = 101995.336
$ws.Range("M63").Value = -101309.336

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H66").Value = 31677.594
$ws.Range("I66").Value = 101995.336
$ws.Range("K66").Value = 509976.68
$ws.Range("M66").Value = -506544.68

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H80").Value = 18585.375
$ws.Range("J80").Value = 18585.375
$ws.Range("L80").Value = 18585.375
$ws.Range("N80").Value = -20581.375

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H83").Value = 18585.375
$ws.Range("J83").Value = 18585.375
$ws.Range("L83").Value = 55756.125
$ws.Range("N83").Value = -65740.125

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H116").Value = 395.54544
$ws.Range("I116").Value = 381.93332
$ws.Range("J116").Value = 424.7143
$ws.Range("K116").Value = 381.93332
$ws.Range("L116").Value = 424.7143
$ws.Range("M116").Value = 1912.06668
$ws.Range("N116").Value = -5012.7143

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H132").Value = 5687.364
$ws.Range("I132").Value = 5998.55
$ws.Range("J132").Value = 5208.615
$ws.Range("K132").Value = 17995.65
$ws.Range("L132").Value = 15625.845
$ws.Range("M132").Value = -15465.65
$ws.Range("N132").Value = -20685.845

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 395.54544
$ws.Range("I3").Value = 381.93332
$ws.Range("J3").Value = 424.7143
$ws.Range("K3").Value = 381.93332
$ws.Range("L3").Value = 424.7143
$ws.Range("M3").Value = -267.93332
$ws.Range("N3").Value = -652.7143

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 2843.1592
$ws.Range("I134").Value = 2809.375
$ws.Range("J134").Value = 2933.25
$ws.Range("K134").Value = 8428.125
$ws.Range("L134").Value = 8799.75
$ws.Range("M134").Value = -5893.125
$ws.Range("N134").Value = -13869.75

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 6068.7573
$ws.Range("I31").Value = 1889.0416
$ws.Range("J31").Value = 8249.478999999999
$ws.Range("K31").Value = 1889.0416
$ws.Range("L31").Value = 8249.478999999999
$ws.Range("M31").Value = -1594.0416
$ws.Range("N31").Value = -8839.478999999999

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 6068.7573
$ws.Range("I34").Value = 1889.0416
$ws.Range("J34").Value = 8249.478999999999
$ws.Range("K34").Value = 1889.0416
$ws.Range("L34").Value = 8249.478999999999
$ws.Range("M34").Value = -1687.0416
$ws.Range("N34").Value = -8653.478999999999

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H62").Value = 4650.8
$ws.Range("I62").Value = 4688.5
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 4688.5
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -4064.5
$ws.Range("N62").Value = -5748

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H65").Value = 4650.8
$ws.Range("I65").Value = 4688.5
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 23442.5
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -20322.5
$ws.Range("N65").Value = -28740

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H97").Value = 29700
$ws.Range("J97").Value = 29650
$ws.Range("L97").Value = 29650
$ws.Range("N97").Value = -31632

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H132").Value = 4387795.5
$ws.Range("I132").Value = 1605.3103
$ws.Range("J132").Value = 18521076
$ws.Range("K132").Value = 4815.9309
$ws.Range("L132").Value = 55563228
$ws.Range("M132").Value = -2285.9309
$ws.Range("N132").Value = -55568288

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H132").Value = 2188.9805
$ws.Range("I132").Value = 1845.4642
$ws.Range("J132").Value = 2607.1738
$ws.Range("K132").Value = 5536.392599999999
$ws.Range("L132").Value = 7821.5214
$ws.Range("M132").Value = -3006.392599999999
$ws.Range("N132").Value = -12881.5214

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H14").Value = 100005
$ws.Range("J14").Value = 100005
$ws.Range("L14").Value = 100005
$ws.Range("N14").Value = -100349

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H132").Value = 2978.5476
$ws.Range("I132").Value = 2581.7407
$ws.Range("J132").Value = 3692.8
$ws.Range("K132").Value = 7745.222099999999
$ws.Range("L132").Value = 11078.4
$ws.Range("M132").Value = -5215.222099999999
$ws.Range("N132").Value = -16138.4

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H14").Value = 502169.7
$ws.Range("I14").Value = 2039.4
$ws.Range("J14").Value = 1002300
$ws.Range("K14").Value = 2039.4
$ws.Range("L14").Value = 1002300
$ws.Range("M14").Value = -1871.4
$ws.Range("N14").Value = -1002636

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H54").Value = 8321.75
$ws.Range("J54").Value = 8929
$ws.Range("L54").Value = 8929
$ws.Range("N54").Value = -9969

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H58").Value = 18000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 18000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 18000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -18616

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H81").Value = 3531.611
$ws.Range("I81").Value = 3376.4285
$ws.Range("J81").Value = 4074.75
$ws.Range("K81").Value = 6752.857
$ws.Range("L81").Value = 8149.5
$ws.Range("M81").Value = -5691.857
$ws.Range("N81").Value = -10271.5

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H84").Value = 3531.611
$ws.Range("I84").Value = 3376.4285
$ws.Range("J84").Value = 4074.75
$ws.Range("K84").Value = 33764.285
$ws.Range("L84").Value = 40747.5
$ws.Range("M84").Value = -28460.285
$ws.Range("N84").Value = -51355.5

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H96").Value = 3944.5454
$ws.Range("I96").Value = 3648.75
$ws.Range("K96").Value = 3648.75
$ws.Range("M96").Value = -2275.75

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H108").Value = 60313
$ws.Range("J108").Value = 60313
$ws.Range("L108").Value = 60313
$ws.Range("N108").Value = -67993
